$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.207.62"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.633.85"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.38"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.256"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.33"
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "1.637.93"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.546"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.25"
$ws.Range("E15").Value = "  -3.51%  "
$ws.Range("D16").Value = "27.170.45"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.03"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.95"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.40"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("E22").Value = "  -6.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.07"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.07"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.68"
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").Value = "1.343.33"
$ws.Range("E33").Value = "  +5.86%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.853"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.73"
$ws.Range("E42").Value = "  +4.41%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.27"
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.773.26"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.86"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("E47").Value = "  +22.12%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0100"
$ws.Range("E48").Value = "  -5.75%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0991"
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.60"
$ws.Range("E51").Value = "  -1.11%  "
